$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "0.9993") are stored as strings, matching the source data which
# uses inline/shared strings for all Price values, not numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.069.82'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '1.789.82'
$ws.Range("E3").Value = '  +2.07%  '
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '327.19'
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").Value = '0.9987'
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("D7").Value = '0.4537'
$ws.Range("E7").Value = '  +1.65%  '
$ws.Range("D8").Value = '0.3602'
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("D9").Value = '0.07517'
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("D10").Value = '42.34'
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("D11").Value = '1.113'
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("D12").Value = '0.9982'
$ws.Range("E12").Value = '  -0.17%  '
$ws.Range("D13").Value = '21.05'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").Value = '6.080'
$ws.Range("E14").Value = '  +0.87%  '
$ws.Range("D15").Value = '7.256'
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = '1.780.92'
$ws.Range("E16").Value = '  +1.63%  '
$ws.Range("D17").Value = '94.13'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '0.00001065'
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = '0.06443'
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").Value = '0.9982'
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("D21").Value = '17.23'
$ws.Range("E21").Value = '  +2.29%  '
$ws.Range("D22").Value = '5.832'
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("D23").Value = '28.086.77'
$ws.Range("E23").Value = '  +1.76%  '
$ws.Range("D24").Value = '11.41'
$ws.Range("E24").Value = '  +1.90%  '
$ws.Range("D25").Value = '2.082'
$ws.Range("E25").Value = '  -1.21%  '
$ws.Range("D26").Value = '163.62'
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("D27").Value = '20.45'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = '1.986.92'
$ws.Range("E28").Value = '  +1.74%  '
$ws.Range("D29").Value = '2.268'
$ws.Range("E29").Value = '  +8.77%  '
$ws.Range("D30").Value = '126.52'
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = '1.116'
$ws.Range("E31").Value = '  +3.51%  '
$ws.Range("D32").Value = '0.09212'
$ws.Range("E32").Value = '  +1.48%  '
$ws.Range("D33").Value = '3.682'
$ws.Range("E33").Value = '  +0.32%  '
$ws.Range("D34").Value = '5.614'
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("D35").Value = '11.97'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = '0.02306'
$ws.Range("E36").Value = '  +0.71%  '
$ws.Range("D37").Value = '0.06156'
$ws.Range("E37").Value = '  +2.03%  '
$ws.Range("D38").Value = '0.2101'
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("D39").Value = '0.6372'
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").Value = '4.998'
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("D41").Value = '1.195'
$ws.Range("E41").Value = '  -0.52%  '
$ws.Range("D42").Value = '1.392'
$ws.Range("E42").Value = '  +1.19%  '
$ws.Range("D43").Value = '7.950'
$ws.Range("E43").Value = '  +2.00%  '
$ws.Range("D44").Value = '13.41'
$ws.Range("E44").Value = '  +1.57%  '
$ws.Range("D45").Value = '0.5944'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").Value = '3.736'
$ws.Range("E46").Value = '  +0.66%  '
$ws.Range("D47").Value = '123.19'
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("D48").Value = '1.972'
$ws.Range("E48").Value = '  +0.95%  '
$ws.Range("D49").Value = '0.06959'
$ws.Range("E49").Value = '  +1.55%  '
$ws.Range("D50").Value = '1.144'
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = '73.11'
$ws.Range("E51").Value = '  +0.88%  '

# Reset the number format back to the default cell style so we don't
# leave a stray "@" text-format style on these cells (the source file
# has no explicit style/format on column D).
$ws.Range("D2:D51").Style = "Normal"
